$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 32, shifting rows 32:65 down to 33:66
# (all existing formatting / column widths / styles on the row follow the insert).
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new weekly record.
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 44587
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112001
$ws.Range("G32").Value = "Berenjena"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 250
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 11000
$ws.Range("M32").Value = 10400
$ws.Range("N32").Value = "$/caja 60 unidades"
$ws.Range("O32").Value = "Provincia de Huasco"
$ws.Range("P32").Value = 173
$ws.Range("Q32").Value = 60
$ws.Range("R32").Value = "Hortaliza"
